# Project Sample Project is saved (SAVE): change cell B11 on the "Rules"
# sheet from the text "R40" to the text "1", keeping its existing
# formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B11")

# Use a scratch cell far away from the used range to stash B11's current
# formatting (borders/fill/font/number format) before we touch anything.
$scratch = $ws.Range("Z1")
$target.Copy()
$scratch.PasteSpecial(-4122)   # xlPasteFormats

# Force the new value to be stored as literal text (not auto-coerced to a
# number) by switching the cell to a text number format before assigning it.
$target.NumberFormat = "@"
$target.Value = "1"

# Restore B11's original formatting (border/fill/font/number format), then
# remove the scratch cell so it leaves no trace in the saved workbook.
$scratch.Copy()
$target.PasteSpecial(-4122)    # xlPasteFormats
$scratch.Clear()

$excel.CutCopyMode = $false
